{"js": "// Prepend \"Design: \" to each answer paragraph (style \"List Bullet\") in the\n// feedback table \u2014 the question cells and header row are left untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"List Bullet\" && para.text && para.text.length > 0) {\n    para.getRange(\"Start\").insertText(\"Design: \", \"Start\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Prepend \"Design: \" to each answer paragraph (style \"List Bullet\") in the\n# feedback table -- the question cells and header row are left untouched.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"List Bullet\" -and $p.Range.Text.Trim().Length -gt 0) {\n        $p.Range.InsertBefore(\"Design: \")\n    }\n}\n"}
